$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Header / source block (rows 1-7 stay in place)
$about.Range("A1").Value = "VoaSL Value of a Statistical Life"
$about.Range("A3").Value = "Source:"
$about.Range("B3").Value = "U.S. Environmental Protection Agency"
$about.Range("B4").Value = 2013
$about.Range("B5").Value = "Frequently Asked Questions on Mortality Risk Valuation"

# B6 becomes the real hyperlink to the EPA FAQ page (anchor #whatvalue).
# Add() sets up the relationship + display text; re-assert the cell text
# afterwards since it is the literal URL (not the " - whatvalue" display).
$about.Hyperlinks.Add($about.Range("B6"), "https://www.epa.gov/environmental-economics/mortality-risk-valuation", "whatvalue", "", "https://www.epa.gov/environmental-economics/mortality-risk-valuation - whatvalue") | Out-Null
$about.Range("B6").Value = "https://www.epa.gov/environmental-economics/mortality-risk-valuation#whatvalue"
$about.Range("B6").Style = "Hiperlink"

# B7 drops its hyperlink styling and becomes the quoted question text
$about.Range("B7").Style = "Normal"
$about.Range("B7").Value = """What value of statistical life does EPA use?"""

# Remove the old row 9 (IPEA pollution-average formula block); everything
# below shifts up by one row.
$about.Rows("9:9").Delete()

# Notes block (now rows 9-13, previously 10-14) keeps its text/styling -
# only the currency-adjustment figures below change.
$about.Range("A16").Value = "We adjust 2006 dollars to 2012 dollars using the following conversion factor:"
$about.Range("A17").Value = 1.141
$about.Range("A18").Value = "See ""cpi.xlsx"" in the InputData folder for source information."

# ---------------------------------------------------------------------
# Sheet "VoaSL"
# ---------------------------------------------------------------------
$voasl = $wb.Worksheets.Item("VoaSL")

$voasl.Range("B2").Style = "Normal"
$voasl.Range("B2").Formula = "=7.4*10^6*About!A17"

# ---------------------------------------------------------------------
# Restore the cursor/selection to what was active when the edit was made
# ---------------------------------------------------------------------
$voasl.Range("A3").Select() | Out-Null
$about.Range("A1").Select() | Out-Null
$about.Activate() | Out-Null
